$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (column D) cells keep their text representation (e.g. trailing
# zeros like "1.000" or thousand-dot formatted values) instead of Excel
# auto-converting the assigned string into a Number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.042.08"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "1.777.59"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "329.45"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.4494"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").Value = "0.3566"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").Value = "0.07455"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").Value = "42.18"
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").Value = "1.108"
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "21.01"
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("D14").Value = "6.061"
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").Value = "7.272"
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("D16").Value = "1.778.95"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").Value = "93.71"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").Value = "0.00001064"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Value = "0.06435"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "17.17"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("D22").Value = "5.805"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").Value = "28.061.93"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").Value = "11.36"
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("D25").Value = "2.125"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("D26").Value = "161.79"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("D27").Value = "20.40"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").Value = "1.983.32"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("D29").Value = "2.164"
$ws.Range("E29").Value = "  +5.79%  "
$ws.Range("D30").Value = "125.05"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").Value = "1.113"
$ws.Range("E31").Value = "  +5.97%  "
$ws.Range("D32").Value = "5.703"
$ws.Range("E32").Value = "  +5.82%  "
$ws.Range("D33").Value = "0.09211"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("D35").Value = "11.90"
$ws.Range("E35").Value = "  +2.37%  "
$ws.Range("D36").Value = "0.06219"
$ws.Range("E36").Value = "  +3.15%  "
$ws.Range("D37").Value = "0.02294"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").Value = "0.2113"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("D39").Value = "5.002"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").Value = "0.6338"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("D41").Value = "1.187"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").Value = "1.398"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("D43").Value = "7.921"
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").Value = "13.35"
$ws.Range("E44").Value = "  +1.83%  "
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").Value = "0.5923"
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("D48").Value = "1.964"
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("D49").Value = "1.144"
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("D50").Value = "0.06901"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").Value = "73.01"
$ws.Range("E51").Value = "  +2.13%  "
